# Leader Board and Analytics added
#
# 1) Text content changes (shared strings used by row 2 of the
#    "E2MTestData" sheet):
#      A2: "Mozila"                    -> "Chrome"
#      H2: "Lacity Talent Games 2018"  -> "Real Me Pre Launch Event"
# 2) Column O (15th column) width tweak: 36.3 -> 36.31
# 3) Sheet view / selection: active cell moves from F2 to I11
#    (and the visible top-left cell shifts from A1 towards D1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the two changed text values -----------------------------
$ws.Range("A2").Value = "Chrome"
$ws.Range("H2").Value = "Real Me Pre Launch Event"

# --- 2) Nudge column O's width (smallest achievable step towards 36.31) -
$ws.Columns.Item(15).ColumnWidth = 35.5

# --- 3) Update the view/selection so the active cell becomes I11 -------
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I11").Select()
